# Updated symbol list on Sat Jan  7 15:09:39 UTC 2023 with GitHub Actions
# Applies updated Price / Volume(1h) / Hora values to the crypto listing sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '260.47'
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '1.58%'
$c.Style = "Normal"
$c = $ws.Range("G2")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '27.25'
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '2.48%'
$c.Style = "Normal"
$c = $ws.Range("G3")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '4.672'
$c.Style = "Normal"
$c = $ws.Range("G4")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '4.09%'
$c.Style = "Normal"
$c = $ws.Range("G5")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '6.664'
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '0.94%'
$c.Style = "Normal"
$c = $ws.Range("G6")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.8500'
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '-0.63%'
$c.Style = "Normal"
$c = $ws.Range("G7")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.9139'
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '0.45%'
$c.Style = "Normal"
$c = $ws.Range("G8")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.1409'
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '2.41%'
$c.Style = "Normal"
$c = $ws.Range("G9")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.04834'
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '8.10%'
$c.Style = "Normal"
$c = $ws.Range("G10")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.07090'
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '1.39%'
$c.Style = "Normal"
$c = $ws.Range("G11")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.03113'
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '3.31%'
$c.Style = "Normal"
$c = $ws.Range("G12")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.09051'
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '-0.64%'
$c.Style = "Normal"
$c = $ws.Range("G13")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.001529'
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '0.27%'
$c.Style = "Normal"
$c = $ws.Range("G14")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.0006162'
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '2.30%'
$c.Style = "Normal"
$c = $ws.Range("G15")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.006028'
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '-0.25%'
$c.Style = "Normal"
$c = $ws.Range("G16")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '3.451'
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '-0.49%'
$c.Style = "Normal"
$c = $ws.Range("G17")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '3.153'
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '0.67%'
$c.Style = "Normal"
$c = $ws.Range("G18")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '1.38%'
$c.Style = "Normal"
$c = $ws.Range("G19")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '0.10%'
$c.Style = "Normal"
$c = $ws.Range("G20")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '0.89%'
$c.Style = "Normal"
$c = $ws.Range("G21")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '4.078'
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '5.91%'
$c.Style = "Normal"
$c = $ws.Range("G22")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.04249'
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '1.26%'
$c.Style = "Normal"
$c = $ws.Range("G23")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '0.04%'
$c.Style = "Normal"
$c = $ws.Range("G24")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '0.003800'
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '-15.04%'
$c.Style = "Normal"
$c = $ws.Range("G25")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '0.10%'
$c.Style = "Normal"
$c = $ws.Range("G26")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '0.0001574'
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '-8.09%'
$c.Style = "Normal"
$c = $ws.Range("G27")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("G28")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("G29")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("G30")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("G31")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("G32")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("G33")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("G34")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("G35")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("G36")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("G37")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("G38")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("G39")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.03871'
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '2.02%'
$c.Style = "Normal"
$c = $ws.Range("G40")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.1113'
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '1.43%'
$c.Style = "Normal"
$c = $ws.Range("G41")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.004079'
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '-34.35%'
$c.Style = "Normal"
$c = $ws.Range("G42")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '8.91%'
$c.Style = "Normal"
$c = $ws.Range("G43")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '2.38%'
$c.Style = "Normal"
$c = $ws.Range("G44")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.00005154'
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '1.42%'
$c.Style = "Normal"
$c = $ws.Range("G45")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '0.10%'
$c.Style = "Normal"
$c = $ws.Range("G46")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '8.13%'
$c.Style = "Normal"
$c = $ws.Range("G47")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("G48")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '0.10%'
$c.Style = "Normal"
$c = $ws.Range("G49")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '0.10%'
$c.Style = "Normal"
$c = $ws.Range("G50")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
$c = $ws.Range("G51")
$c.NumberFormat = "@"
$c.Value = '15'
$c.Style = "Normal"
